# Fix thumbnail URL generation in catalog: each row's ThumbnailURL (column I)
# should point to the chart image for that row's own county/region (column E),
# not always to Delaware's chart. Column I text is corrected in-place; the row
# order (Delaware rows already correct, others grouped by chart-type block)
# mirrors the original authoring order so new shared strings are appended in
# the same sequence as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$base = "https://raw.githubusercontent.com/morpc-insights/insights-crashes/refs/heads/main/output_data/charts/"

# County/region order as laid out in the sheet (row offset 0 = Delaware, already correct)
$counties = @("Delaware", "Fairfield", "Franklin", "Knox", "Licking", "Madison", "Marion", "Morrow", "Pickaway", "Union", "Region10")

# Each block: chart-type slug used in the filename, and the first row of the block
$blocks = @(
    @{ Slug = "crashSeverityByType"; StartRow = 2 },
    @{ Slug = "crashSeverityByYear"; StartRow = 13 },
    @{ Slug = "crashSeverityByMode"; StartRow = 24 }
)

foreach ($block in $blocks) {
    for ($i = 1; $i -lt $counties.Count; $i++) {
        $row = $block.StartRow + $i
        $county = $counties[$i]
        $url = $base + $block.Slug + "-" + $county + ".svg"
        $ws.Range("I" + $row).Value = $url
    }
}
